$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Learning Path")

# Insert a new header row at the top, pushing the two existing rows down by one
$ws.Rows.Item(1).Insert()

# Header row: Item / Date / Status, bold
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Status"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)

# New "UFO dataset from Acloud guru" tracker rows
$ws.Range("A4").Value = "Complete credit-fraud-dealing-with-imbalanced-datasets - https://www.kaggle.com/janiobachmann/credit-fraud-dealing-with-imbalanced-datasets"
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 43.5

$ws.Range("B4").Value = 43790
$ws.Range("B4").NumberFormat = "mm-dd-yy"
$ws.Range("C4").Value = "In Progress"

$ws.Range("A5").Value = "Do a spark/Datbricks Project"
$ws.Range("A6").Value = "Understand PCA blood and bone means corrctly"

# Propagate the date format to B5:B6 and status text to C5:C6 by copying the
# already-styled cells, so the same style record is reused instead of new
# duplicate ones being allocated.
$ws.Range("B4").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)
$ws.Range("B5").Value = 43790
$ws.Range("B6").Value = 43790

$ws.Range("C4").Copy()
$ws.Range("C5:C6").PasteSpecial(-4122)
$ws.Range("C5").Value = "In Progress"
$ws.Range("C6").Value = "In Progress"

$ws.Range("A7").Select() | Out-Null
